$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add D1 "rename" with header formatting copied from C1 (bold/border style)
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "rename"

# Copy the bold/bordered "index" style (used in column A) as a template for new A-column cells
$ws.Range("A2").Copy() | Out-Null

$ws.Range("A2").Value = 2
$ws.Range("B2").Value = "coulomb"
$ws.Range("C2").Value = 0.0003168259644355218
$ws.Range("D2").Value = "`$\langle pp \vert \vert qq \rangle`$"

$ws.Range("A3").Value = 3
$ws.Range("B3").Value = "screen1_1"
$ws.Range("C3").Value = 0.00005940879602786336
$ws.Range("D3").Value = "`$(\langle pp \vert \vert rr \rangle)_{1}`$"

$ws.Range("A4").Value = 4
$ws.Range("B4").Value = "screen1_2"
$ws.Range("C4").Value = 0.00007824272819801549
$ws.Range("D4").Value = "`$(\langle pp \vert \vert rr \rangle)_{2}`$"

$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "screen1_3"
$ws.Range("C5").Value = 0.00006253773474039816
$ws.Range("D5").Value = "`$(\langle pp \vert \vert rr \rangle)_{3}`$"

$ws.Range("A6").Value = 6
$ws.Range("B6").Value = "screen1_4"
$ws.Range("C6").Value = 0.0001148897168160615
$ws.Range("D6").Value = "`$(\langle pp \vert \vert rr \rangle)_{4}`$"

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = "screen2_1"
$ws.Range("C7").Value = 0.0001258499164456318
$ws.Range("D7").Value = "`$(\langle qq \vert \vert ss \rangle)_{1}`$"

$ws.Range("A8").Value = 8
$ws.Range("B8").Value = "screen2_2"
$ws.Range("C8").Value = 0.0000735040271998253
$ws.Range("D8").Value = "`$(\langle qq \vert \vert ss \rangle)_{2}`$"

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "screen2_3"
$ws.Range("C9").Value = 0.00004517908053822079
$ws.Range("D9").Value = "`$(\langle qq \vert \vert ss \rangle)_{3}`$"

$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "screen2_4"
$ws.Range("C10").Value = 0.00003661295765833431
$ws.Range("D10").Value = "`$(\langle qq \vert \vert ss \rangle)_{4}`$"

$ws.Range("A11").Value = 11
$ws.Range("B11").Value = "eijab_1"
$ws.Range("C11").Value = 0.00003889256167491719
$ws.Range("D11").Value = "`$(e^{rs}_{pq})_{1}`$"

$ws.Range("A12").Value = 12
$ws.Range("B12").Value = "eijab_2"
$ws.Range("C12").Value = 0.00003276834448683038
$ws.Range("D12").Value = "`$(e^{rs}_{pq})_{2}`$"

$ws.Range("A13").Value = 13
$ws.Range("B13").Value = "eijab_3"
$ws.Range("C13").Value = 0.00001989628732556284
$ws.Range("D13").Value = "`$(e^{rs}_{pq})_{3}`$"

$ws.Range("A14").Value = 14
$ws.Range("B14").Value = "eijab_4"
$ws.Range("C14").Value = 0.00002793296040015019
$ws.Range("D14").Value = "`$(e^{rs}_{pq})_{4}`$"

$ws.Range("A15").Value = 15
$ws.Range("B15").Value = "screenvirt_1"
$ws.Range("C15").Value = 0.00003905361518297251
$ws.Range("D15").Value = "`$(\langle ss \vert \vert rr \rangle)_{1}`$"

$ws.Range("A16").Value = 16
$ws.Range("B16").Value = "screenvirt_2"
$ws.Range("C16").Value = 0.0000378460080171781
$ws.Range("D16").Value = "`$(\langle ss \vert \vert rr \rangle)_{2}`$"

$ws.Range("A17").Value = 17
$ws.Range("B17").Value = "screenvirt_3"
$ws.Range("C17").Value = 0.00007208968252806431
$ws.Range("D17").Value = "`$(\langle ss \vert \vert rr \rangle)_{3}`$"

$ws.Range("A18").Value = 18
$ws.Range("B18").Value = "screenvirt_4"
$ws.Range("C18").Value = 0.00006732142213003006
$ws.Range("D18").Value = "`$(\langle ss \vert \vert rr \rangle)_{4}`$"

$ws.Range("A19").Value = 19
$ws.Range("B19").Value = "Fr1"
$ws.Range("C19").Value = 0.00003859239089305153
$ws.Range("D19").Value = "`$(F_{r})_{1}`$"

$ws.Range("A20").Value = 20
$ws.Range("B20").Value = "Fr2"
$ws.Range("C20").Value = 0.0000227742806884989
$ws.Range("D20").Value = "`$(F_{r})_{2}`$"

$ws.Range("A21").Value = 22
$ws.Range("B21").Value = "Fr4"
$ws.Range("C21").Value = 0.00001634412636579059
$ws.Range("D21").Value = "`$(F_{r})_{4}`$"

$ws.Range("A22").Value = 23
$ws.Range("B22").Value = "Fs1"
$ws.Range("C22").Value = 0.0001776322197255775
$ws.Range("D22").Value = "`$(F_{s})_{1}`$"

$ws.Range("A23").Value = 25
$ws.Range("B23").Value = "Fs3"
$ws.Range("C23").Value = 0.00002693023691802636
$ws.Range("D23").Value = "`$(F_{s})_{3}`$"

$ws.Range("A24").Value = 26
$ws.Range("B24").Value = "Fs4"
$ws.Range("C24").Value = 0.00005324781947464344
$ws.Range("D24").Value = "`$(F_{s})_{4}`$"

$ws.Range("A25").Value = 27
$ws.Range("B25").Value = "occr1"
$ws.Range("C25").Value = 0.00004265310876794375
$ws.Range("D25").Value = "`$(\eta_{r})_{1}`$"

$ws.Range("A26").Value = 28
$ws.Range("B26").Value = "occr2"
$ws.Range("C26").Value = 0.00004020503578892763
$ws.Range("D26").Value = "`$(\eta_{r})_{2}`$"

$ws.Range("A27").Value = 29
$ws.Range("B27").Value = "occr3"
$ws.Range("C27").Value = 0.00002287255070539599
$ws.Range("D27").Value = "`$(\eta_{r})_{3}`$"

$ws.Range("A28").Value = 30
$ws.Range("B28").Value = "occr4"
$ws.Range("C28").Value = 0.00002912423636206734
$ws.Range("D28").Value = "`$(\eta_{r})_{4}`$"

$ws.Range("A29").Value = 33
$ws.Range("B29").Value = "occs3"
$ws.Range("C29").Value = 0.00001656437647339937
$ws.Range("D29").Value = "`$(\eta_{s})_{3}`$"

$ws.Range("A30").Value = 34
$ws.Range("B30").Value = "occs4"
$ws.Range("C30").Value = 0.00001621568835734885
$ws.Range("D30").Value = "`$(\eta_{s})_{4}`$"

$ws.Range("A31").Value = 35
$ws.Range("B31").Value = "SCFFr1"
$ws.Range("C31").Value = 0.00001633159600448834
$ws.Range("D31").Value = "`$(F_{r}^{\text{SCF}})_{1}`$"

$ws.Range("A32").Value = 37
$ws.Range("B32").Value = "SCFFr3"
$ws.Range("C32").Value = 0.00001133333531881307
$ws.Range("D32").Value = "`$(F_{r}^{\text{SCF}})_{3}`$"

$ws.Range("A33").Value = 39
$ws.Range("B33").Value = "SCFFs1"
$ws.Range("C33").Value = 0.0002830631350236944
$ws.Range("D33").Value = "`$(F_{s}^{\text{SCF}})_{1}`$"

$ws.Range("A34").Value = 40
$ws.Range("B34").Value = "SCFFs2"
$ws.Range("C34").Value = 0.0001792907481552214
$ws.Range("D34").Value = "`$(F_{s}^{\text{SCF}})_{2}`$"

$ws.Range("A35").Value = 41
$ws.Range("B35").Value = "SCFFs3"
$ws.Range("C35").Value = 0.00001446614255483271
$ws.Range("D35").Value = "`$(F_{s}^{\text{SCF}})_{3}`$"

$ws.Range("A36").Value = 42
$ws.Range("B36").Value = "SCFFs4"
$ws.Range("C36").Value = 0.0000298289304025045
$ws.Range("D36").Value = "`$(F_{s}^{\text{SCF}})_{4}`$"

$ws.Range("A37").Value = 51
$ws.Range("B37").Value = "hrr1"
$ws.Range("C37").Value = 0.00001896106418495491
$ws.Range("D37").Value = "`$(h_{rr})_{1}`$"

$ws.Range("A38").Value = 52
$ws.Range("B38").Value = "hrr2"
$ws.Range("C38").Value = 0.00001151941430186525
$ws.Range("D38").Value = "`$(h_{rr})_{2}`$"

$ws.Range("A39").Value = 53
$ws.Range("B39").Value = "hrr3"
$ws.Range("C39").Value = 0.00001809523655504705
$ws.Range("D39").Value = "`$(h_{rr})_{3}`$"

$ws.Range("A40").Value = 55
$ws.Range("B40").Value = "hss1"
$ws.Range("C40").Value = 0.00003099170070176877
$ws.Range("D40").Value = "`$(h_{ss})_{1}`$"

$ws.Range("A41").Value = 56
$ws.Range("B41").Value = "hss2"
$ws.Range("C41").Value = 0.0001180824240041617
$ws.Range("D41").Value = "`$(h_{ss})_{2}`$"

$ws.Range("A42").PasteSpecial(-4122) | Out-Null
$ws.Range("A42").Value = 59
$ws.Range("B42").Value = "hpp"
$ws.Range("C42").Value = 0.0001177320874143567
$ws.Range("D42").Value = "`$h_{pp}`$"

$ws.Range("A43").PasteSpecial(-4122) | Out-Null
$ws.Range("A43").Value = 61
$ws.Range("B43").Value = "Fp"
$ws.Range("C43").Value = 0.00004691652598230687
$ws.Range("D43").Value = "`$F_{p}`$"

$ws.Range("A44").PasteSpecial(-4122) | Out-Null
$ws.Range("A44").Value = 62
$ws.Range("B44").Value = "Fq"
$ws.Range("C44").Value = 0.00005396310224662725
$ws.Range("D44").Value = "`$F_{q}`$"

$ws.Range("A45").PasteSpecial(-4122) | Out-Null
$ws.Range("A45").Value = 64
$ws.Range("B45").Value = "occq"
$ws.Range("C45").Value = 0.00009934645994097476
$ws.Range("D45").Value = "`$\eta_{q}`$"

$ws.Range("A46").PasteSpecial(-4122) | Out-Null
$ws.Range("A46").Value = 65
$ws.Range("B46").Value = "SCFFp"
$ws.Range("C46").Value = 0.0001755466996723631
$ws.Range("D46").Value = "`$F_{p}^{\text{SCF}}`$"

$ws.Range("A47").PasteSpecial(-4122) | Out-Null
$ws.Range("A47").Value = 66
$ws.Range("B47").Value = "SCFFq"
$ws.Range("C47").Value = 0.00004472042494480772
$ws.Range("D47").Value = "`$F_{q}^{\text{SCF}}`$"

Write-Output "done"